$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column H (Genero) / I (Estado) so that:
#   old H (Genero) -> J, old I (Estado) -> K
# This mirrors Excel's native "Insert Columns" behaviour (keeps the
# 25.28515625-wide column formatting continuity through column J).
$ws.Columns.Item(8).Resize(1, 2).EntireColumn.Insert()

# Fill column H (Nacionalidad) top-to-bottom first ...
$ws.Range("H1").Value = "Nacionalidad"
$ws.Range("H2").Value = "Colombian"
$ws.Range("H3").Value = "Angolan"
$ws.Range("H4").Value = "Bolivian"
$ws.Range("H5").Value = "Peruvian"
$ws.Range("H6").Value = "Spanish"
$ws.Range("H7").Value = "American"

# ... then column I (Estado Civil) top-to-bottom
$ws.Range("I1").Value = "Estado Civil"
$ws.Range("I2").Value = "Single"
$ws.Range("I3").Value = "Married"
$ws.Range("I4").Value = "Other"
$ws.Range("I5").Value = "Married"
$ws.Range("I6").Value = "Single"
$ws.Range("I7").Value = "Married"

# Restore the quote-prefix style (s="1") that Excel's own column-insert
# carries into newly inserted cells, matching the rest of the data rows'
# formatting (lost because setting .Value resets the cell style here).
$ws.Range("F2:G7").Copy()
$ws.Range("H2:I7").PasteSpecial(-4122)

# Data correction: row 3 (Juan Carlos Diaz Perez) Estado changed to "Usado"
$ws.Range("K3").Value = "Usado"

# Update the active selection to reflect the last-edited cell
$ws.Range("I7").Select()
